$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then column letter -> new text value
$updates = @(
    @{Row=2; D='60.153.17'; E='  +1.11%  '},
    @{Row=3; D='2.594.65'; E='  +0.23%  '},
    @{Row=4; E='  +0.07%  '},
    @{Row=5; D='578.53'; E='  +4.42%  '},
    @{Row=6; D='142.62'; E='  +1.47%  '},
    @{Row=7; D='0.999'; E='  +0.12%  '},
    @{Row=8; D='0.596'; E='  +0.52%  '},
    @{Row=9; D='2.600.76'; E='  -0.17%  '},
    @{Row=10; E='  -2.98%  '},
    @{Row=11; D='0.105'; E='  +1.11%  '},
    @{Row=12; E='  -2.65%  '},
    @{Row=13; D='0.370'; E='  +3.65%  '},
    @{Row=14; D='3.055.89'; E='  +0.26%  '},
    @{Row=15; D='24.67'; E='  +7.16%  '},
    @{Row=16; D='60.176.73'; E='  +1.19%  '},
    @{Row=17; E='  +2.69%  '},
    @{Row=18; D='2.601.52'; E='  +0.34%  '},
    @{Row=19; D='11.51'; E='  +10.47%  '},
    @{Row=20; D='4.64'; E='  +1.70%  '},
    @{Row=21; D='346.54'; E='  +1.90%  '},
    @{Row=22; D='6.89'; E='  +4.73%  '},
    @{Row=23; D='0.999'; E='  -0.10%  '},
    @{Row=24; D='0.524'; E='  +8.57%  '},
    @{Row=25; D='63.01'; E='  +0.03%  '},
    @{Row=26; E='  +0.40%  '},
    @{Row=27; E='  +0.16%  '},
    @{Row=28; D='8.04'; E='  +7.53%  '},
    @{Row=29; D='0.0₃0791'; E='  +2.74%  '},
    @{Row=30; E='  +10.80%  '},
    @{Row=31; B='USDe'; C='https://coinranking.com/coin/exbfr2U-0+usde-usde'; D='0.999'; E='  +0.08%  '},
    @{Row=32; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='6.36'; E='  +4.14%  '},
    @{Row=33; D='163.58'; E='  +3.80%  '},
    @{Row=34; D='19.40'; E='  +0.17%  '},
    @{Row=35; D='4.29'; E='  +5.53%  '},
    @{Row=36; D='0.988'; E='  +7.17%  '},
    @{Row=37; D='1.24'; E='  +7.08%  '},
    @{Row=38; D='1.62'; E='  +9.15%  '},
    @{Row=39; D='37.98'; E='  +1.07%  '},
    @{Row=40; D='3.90'; E='  +5.97%  '},
    @{Row=41; D='308.02'; E='  +6.37%  '},
    @{Row=42; D='0.836'; E='  -0.33%  '},
    @{Row=43; D='135.42'; E='  -0.79%  '},
    @{Row=44; E='  +0.24%  '},
    @{Row=45; D='0.0988'; E='  +1.38%  '},
    @{Row=46; D='5.01'; E='  +10.46%  '},
    @{Row=47; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='19.70'; E='  +3.92%  '},
    @{Row=48; B='Mantle'; C='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D='0.602'; E='  +0.62%  '},
    @{Row=49; D='0.0548'; E='  +2.45%  '},
    @{Row=50; D='20.06'; E='  +7.82%  '},
    @{Row=51; E='  +2.13%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey('D')) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }
}
